$d = $word.ActiveDocument
$d.Content.Find.Execute("Week 2:", $true, $false, $false, $false, $false, $true, 1, $false, "Week 4:", 2)
